$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 153.3077
$ws.Range("I2").Value = 166.45454
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 166.45454
$ws.Range("L2").Value = 81
$ws.Range("M2").Value = -53.45454000000001
$ws.Range("N2").Value = -307
$ws.Range("H19").Value = 586.625
$ws.Range("J19").Value = 478.4
$ws.Range("L19").Value = 478.4
$ws.Range("N19").Value = -828.4
$ws.Range("H33").Value = 423.42856
$ws.Range("J33").Value = 542
$ws.Range("L33").Value = 542
$ws.Range("N33").Value = -1000
$ws.Range("H112").Value = 3911.111
$ws.Range("J112").Value = 4714.2856
$ws.Range("L112").Value = 14142.8568
$ws.Range("N112").Value = -16358.8568
$ws.Range("H116").Value = 3097
$ws.Range("I116").Value = 2573.6365
$ws.Range("J116").Value = 3919.4285
$ws.Range("K116").Value = 2573.6365
$ws.Range("L116").Value = 3919.4285
$ws.Range("M116").Value = 868.3634999999999
$ws.Range("N116").Value = -10803.4285
$ws.Range("H137").Value = 1466.2433
$ws.Range("I137").Value = 1015.5
$ws.Range("K137").Value = 3046.5
$ws.Range("M137").Value = -496.5
$ws.Range("H138").Value = 1582.32
$ws.Range("J138").Value = 1820.5385
$ws.Range("L138").Value = 5461.6155
$ws.Range("N138").Value = -15741.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 686.375
$ws.Range("I2").Value = 607.15
$ws.Range("K2").Value = 607.15
$ws.Range("M2").Value = -494.15
$ws.Range("H32").Value = 4421.3477
$ws.Range("I32").Value = 4334.4546
$ws.Range("J32").Value = 6333
$ws.Range("K32").Value = 4334.4546
$ws.Range("L32").Value = 6333
$ws.Range("M32").Value = -4047.4546
$ws.Range("N32").Value = -6907
$ws.Range("H45").Value = 1107.125
$ws.Range("I45").Value = 1047.2307
$ws.Range("K45").Value = 1047.2307
$ws.Range("M45").Value = -670.2307000000001
$ws.Range("H53").Value = 1012
$ws.Range("J53").Value = 1012
$ws.Range("L53").Value = 1012
$ws.Range("N53").Value = -2376
$ws.Range("H74").Value = 2328.7
$ws.Range("I74").Value = 1328.1428
$ws.Range("J74").Value = 4663.3335
$ws.Range("K74").Value = 1328.1428
$ws.Range("L74").Value = 4663.3335
$ws.Range("M74").Value = -454.1428000000001
$ws.Range("N74").Value = -6411.3335
$ws.Range("H77").Value = 2328.7
$ws.Range("I77").Value = 1328.1428
$ws.Range("J77").Value = 4663.3335
$ws.Range("K77").Value = 6640.714
$ws.Range("L77").Value = 23316.6675
$ws.Range("M77").Value = -2272.714
$ws.Range("N77").Value = -32052.6675
$ws.Range("H116").Value = 686.375
$ws.Range("I116").Value = 607.15
$ws.Range("K116").Value = 607.15
$ws.Range("M116").Value = 1686.85
$ws.Range("H122").Value = 1728.9736
$ws.Range("I122").Value = 1695.9615
$ws.Range("K122").Value = 5087.8845
$ws.Range("M122").Value = -2637.8845
$ws.Range("H132").Value = 2882.2173
$ws.Range("I132").Value = 2427.353
$ws.Range("K132").Value = 7282.059
$ws.Range("M132").Value = -4752.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 686.375
$ws.Range("I3").Value = 607.15
$ws.Range("K3").Value = 607.15
$ws.Range("M3").Value = -493.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1555.8462
$ws.Range("I31").Value = 1518.8334
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1518.8334
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -1223.8334
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 1555.8462
$ws.Range("I34").Value = 1518.8334
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1518.8334
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -1316.8334
$ws.Range("N34").Value = -2404
$ws.Range("H122").Value = 5344.609
$ws.Range("I122").Value = 5496.636
$ws.Range("K122").Value = 16489.908
$ws.Range("M122").Value = -14039.908
$ws.Range("H132").Value = 2105.9333
$ws.Range("I132").Value = 1243.4445
$ws.Range("J132").Value = 3399.6667
$ws.Range("K132").Value = 3730.3335
$ws.Range("L132").Value = 10199.0001
$ws.Range("M132").Value = -1200.3335
$ws.Range("N132").Value = -15259.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 912.4167
$ws.Range("J2").Value = 2122
$ws.Range("L2").Value = 12732
$ws.Range("N2").Value = -12958
$ws.Range("H113").Value = 654.0968
$ws.Range("I113").Value = 585.5
$ws.Range("J113").Value = 677.95654
$ws.Range("K113").Value = 1756.5
$ws.Range("L113").Value = 2033.86962
$ws.Range("M113").Value = 413.5
$ws.Range("N113").Value = -6373.869619999999
$ws.Range("H131").Value = 10990083
$ws.Range("I131").Value = 166667040
$ws.Range("J131").Value = 1120.9529
$ws.Range("K131").Value = 500001120
$ws.Range("L131").Value = 3362.8587
$ws.Range("M131").Value = -499996080
$ws.Range("N131").Value = -13442.8587
$ws.Range("H132").Value = 1016.0833
$ws.Range("I132").Value = 969.3
$ws.Range("K132").Value = 8723.699999999999
$ws.Range("M132").Value = -6193.699999999999
$ws.Range("H140").Value = 24704.936
$ws.Range("I140").Value = 52272.8
$ws.Range("J140").Value = 3498.8845
$ws.Range("K140").Value = 156818.4
$ws.Range("L140").Value = 10496.6535
$ws.Range("M140").Value = -151638.4
$ws.Range("N140").Value = -20856.6535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18003630
$ws.Range("I70").Value = 17860776
$ws.Range("J70").Value = 18185446
$ws.Range("K70").Value = 17860776
$ws.Range("L70").Value = 18185446
$ws.Range("M70").Value = -17860506
$ws.Range("N70").Value = -18185986
$ws.Range("H73").Value = 18003630
$ws.Range("I73").Value = 17860776
$ws.Range("J73").Value = 18185446
$ws.Range("K73").Value = 17860776
$ws.Range("L73").Value = 18185446
$ws.Range("M73").Value = -17859840
$ws.Range("N73").Value = -18187318
$ws.Range("H132").Value = 2677.6562
$ws.Range("I132").Value = 2513
$ws.Range("J132").Value = 2992
$ws.Range("K132").Value = 7539
$ws.Range("L132").Value = 8976
$ws.Range("M132").Value = -5009
$ws.Range("N132").Value = -14036

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3062.6428
$ws.Range("I132").Value = 3261.125
$ws.Range("K132").Value = 9783.375
$ws.Range("M132").Value = -7253.375
